$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1171850.9
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 1171850.9
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 3515552.7
$ws.Range("M17").ClearContents()
$ws.Range("N17").Value = -3515888.7
$ws.Range("H43").Value = 3725
$ws.Range("J43").Value = 2950
$ws.Range("L43").Value = 2950
$ws.Range("N43").Value = -3088
$ws.Range("H55").Value = 197.1
$ws.Range("I55").Value = 146.33333
$ws.Range("J55").Value = 273.25
$ws.Range("K55").Value = 146.33333
$ws.Range("L55").Value = 273.25
$ws.Range("M55").Value = 67.66667000000001
$ws.Range("N55").Value = -701.25
$ws.Range("H98").Value = 1248.1305
$ws.Range("I98").Value = 1285.4
$ws.Range("J98").Value = 999.6667
$ws.Range("K98").Value = 1285.4
$ws.Range("L98").Value = 999.6667
$ws.Range("M98").Value = 212.5999999999999
$ws.Range("N98").Value = -3995.6667
$ws.Range("H106").Value = 1212.5333
$ws.Range("I106").Value = 1245.5714
$ws.Range("K106").Value = 1245.5714
$ws.Range("M106").Value = -614.5714
$ws.Range("H113").Value = 4009.5557
$ws.Range("I113").Value = 1900.6666
$ws.Range("J113").Value = 5064
$ws.Range("K113").Value = 1900.6666
$ws.Range("L113").Value = 5064
$ws.Range("M113").Value = 1353.3334
$ws.Range("N113").Value = -11572
$ws.Range("H122").Value = 1248.1305
$ws.Range("I122").Value = 1285.4
$ws.Range("J122").Value = 999.6667
$ws.Range("K122").Value = 3856.2
$ws.Range("L122").Value = 2999.0001
$ws.Range("M122").Value = -1406.2
$ws.Range("N122").Value = -7899.0001
$ws.Range("H125").Value = 2064.0557
$ws.Range("I125").Value = 1555.4
$ws.Range("K125").Value = 13998.6
$ws.Range("M125").Value = -11538.6
$ws.Range("H132").Value = 5176.154
$ws.Range("I132").Value = 5700.6665
$ws.Range("K132").Value = 17101.9995
$ws.Range("M132").Value = -14571.9995
$ws.Range("H138").Value = 7411737
$ws.Range("J138").Value = 12352149
$ws.Range("L138").Value = 37056447
$ws.Range("N138").Value = -37066727

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1002
$ws.Range("I2").Value = 1002
$ws.Range("K2").Value = 1002
$ws.Range("M2").Value = -889
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 0
$ws.Range("K26").Value = 0
$ws.Range("M26").ClearContents()
$ws.Range("H32").Value = 21127.883
$ws.Range("I32").Value = 19921.182
$ws.Range("K32").Value = 19921.182
$ws.Range("M32").Value = -19634.182
$ws.Range("H45").Value = 1999
$ws.Range("I45").Value = 1790.9
$ws.Range("K45").Value = 1790.9
$ws.Range("M45").Value = -1413.9
$ws.Range("H116").Value = 1002
$ws.Range("I116").Value = 1002
$ws.Range("K116").Value = 1002
$ws.Range("M116").Value = 1292
$ws.Range("H119").Value = 59499.5
$ws.Range("J119").Value = 59499.5
$ws.Range("L119").Value = 59499.5
$ws.Range("N119").Value = -69175.5

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1002
$ws.Range("I3").Value = 1002
$ws.Range("K3").Value = 1002
$ws.Range("M3").Value = -888
$ws.Range("H20").Value = 2174.3076
$ws.Range("I20").Value = 2252
$ws.Range("K20").Value = 2252
$ws.Range("M20").Value = -2005
$ws.Range("H105").Value = 15586.429
$ws.Range("I105").Value = 15586.429
$ws.Range("K105").Value = 15586.429
$ws.Range("M105").Value = -13839.429
$ws.Range("H107").Value = 4777
$ws.Range("I107").Value = 4749.25
$ws.Range("K107").Value = 4749.25
$ws.Range("M107").Value = -2829.25

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4378.483
$ws.Range("I31").Value = 3607.9167
$ws.Range("J31").Value = 4922.4116
$ws.Range("K31").Value = 3607.9167
$ws.Range("L31").Value = 4922.4116
$ws.Range("M31").Value = -3312.9167
$ws.Range("N31").Value = -5512.4116
$ws.Range("H34").Value = 4378.483
$ws.Range("I34").Value = 3607.9167
$ws.Range("J34").Value = 4922.4116
$ws.Range("K34").Value = 3607.9167
$ws.Range("L34").Value = 4922.4116
$ws.Range("M34").Value = -3405.9167
$ws.Range("N34").Value = -5326.4116
$ws.Range("H41").Value = 23755.5
$ws.Range("I41").Value = 5500
$ws.Range("J41").Value = 42011
$ws.Range("K41").Value = 5500
$ws.Range("L41").Value = 42011
$ws.Range("M41").Value = -5072
$ws.Range("N41").Value = -42867
$ws.Range("H47").Value = 10000
$ws.Range("I47").Value = 10000
$ws.Range("K47").Value = 10000
$ws.Range("M47").Value = -9434
$ws.Range("H58").Value = 3605.389
$ws.Range("I58").Value = 3576.7334
$ws.Range("J58").Value = 3748.6667
$ws.Range("K58").Value = 3576.7334
$ws.Range("L58").Value = 3748.6667
$ws.Range("M58").Value = -3373.7334
$ws.Range("N58").Value = -4154.6667
$ws.Range("H122").Value = 2521
$ws.Range("I122").Value = 2194.6667
$ws.Range("K122").Value = 6584.000100000001
$ws.Range("M122").Value = -4134.000100000001
$ws.Range("H136").Value = 3605.389
$ws.Range("I136").Value = 3576.7334
$ws.Range("J136").Value = 3748.6667
$ws.Range("K136").Value = 10730.2002
$ws.Range("L136").Value = 11246.0001
$ws.Range("M136").Value = -8180.200199999999
$ws.Range("N136").Value = -16346.0001
$ws.Range("H141").Value = 197068.08
$ws.Range("J141").Value = 197068.08
$ws.Range("L141").Value = 197068.08
$ws.Range("N141").Value = -207428.08

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1002
$ws.Range("I5").Value = 1012.4
$ws.Range("J5").Value = 950
$ws.Range("K5").Value = 3037.2
$ws.Range("L5").Value = 2850
$ws.Range("M5").Value = -2925.2
$ws.Range("N5").Value = -3074
$ws.Range("H96").Value = 5996.3335
$ws.Range("J96").Value = 5996.3335
$ws.Range("L96").Value = 17989.0005
$ws.Range("N96").Value = -22107.0005
$ws.Range("H97").Value = 262.25
$ws.Range("J97").Value = 249.5
$ws.Range("L97").Value = 748.5
$ws.Range("N97").Value = -1740.5
$ws.Range("H100").Value = 5999.5
$ws.Range("J100").Value = 5999.5
$ws.Range("L100").Value = 17998.5
$ws.Range("N100").Value = -19620.5
$ws.Range("H102").Value = 5000
$ws.Range("J102").Value = 5000
$ws.Range("L102").Value = 15000
$ws.Range("N102").Value = -19868
$ws.Range("H119").Value = 0
$ws.Range("I119").Value = 0
$ws.Range("K119").Value = 0
$ws.Range("M119").ClearContents()
$ws.Range("H120").Value = 8860.714
$ws.Range("I120").Value = 6405.4
$ws.Range("K120").Value = 19216.2
$ws.Range("M120").Value = -14378.2
$ws.Range("H121").Value = 797.1429000000001
$ws.Range("I121").Value = 407.5
$ws.Range("J121").Value = 1316.6666
$ws.Range("K121").Value = 1222.5
$ws.Range("L121").Value = 3949.9998
$ws.Range("M121").Value = 87.5
$ws.Range("N121").Value = -6569.9998
$ws.Range("H131").Value = 1615.1621
$ws.Range("I131").Value = 1183.3334
$ws.Range("K131").Value = 3550.0002
$ws.Range("M131").Value = 1489.9998
$ws.Range("H135").Value = 1002
$ws.Range("I135").Value = 1012.4
$ws.Range("J135").Value = 950
$ws.Range("K135").Value = 9111.6
$ws.Range("L135").Value = 8550
$ws.Range("M135").Value = -6576.6
$ws.Range("N135").Value = -13620

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("N15").ClearContents()
$ws.Range("H70").Value = 403403.2
$ws.Range("H73").Value = 403403.2
$ws.Range("H80").Value = 3234.9092
$ws.Range("I80").Value = 3085
$ws.Range("J80").Value = 3414.8
$ws.Range("K80").Value = 3085
$ws.Range("L80").Value = 3414.8
$ws.Range("M80").Value = -2087
$ws.Range("N80").Value = -5410.8
$ws.Range("H81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").ClearContents()
$ws.Range("H83").Value = 3234.9092
$ws.Range("I83").Value = 3085
$ws.Range("J83").Value = 3414.8
$ws.Range("K83").Value = 15425
$ws.Range("L83").Value = 17074
$ws.Range("M83").Value = -10433
$ws.Range("N83").Value = -27058
$ws.Range("H84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").ClearContents()
$ws.Range("H122").Value = 3561.3635
$ws.Range("I122").Value = 3253.5715
$ws.Range("J122").Value = 4100
$ws.Range("K122").Value = 9760.7145
$ws.Range("L122").Value = 12300
$ws.Range("M122").Value = -7310.7145
$ws.Range("N122").Value = -17200

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3377.7778
$ws.Range("I40").Value = 2914.2856
$ws.Range("K40").Value = 2914.2856
$ws.Range("M40").Value = -2778.2856
$ws.Range("H68").Value = 5894.5
$ws.Range("I68").Value = 2790
$ws.Range("J68").Value = 8999
$ws.Range("K68").Value = 2790
$ws.Range("L68").Value = 8999
$ws.Range("M68").Value = -2041
$ws.Range("N68").Value = -10497
$ws.Range("H71").Value = 5894.5
$ws.Range("I71").Value = 2790
$ws.Range("J71").Value = 8999
$ws.Range("K71").Value = 13950
$ws.Range("L71").Value = 44995
$ws.Range("M71").Value = -10206
$ws.Range("N71").Value = -52483
$ws.Range("H129").Value = 50000
$ws.Range("J129").Value = 50000
$ws.Range("L129").Value = 50000
$ws.Range("N129").Value = -60000

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 341.57144
$ws.Range("J107").Value = 0
$ws.Range("L107").Value = 0
$ws.Range("N107").ClearContents()
$ws.Range("H136").Value = 6227
$ws.Range("I136").Value = 1704
$ws.Range("J136").Value = 7131.6
$ws.Range("K136").Value = 5112
$ws.Range("L136").Value = 21394.8
$ws.Range("M136").Value = -2562
$ws.Range("N136").Value = -26494.8
